$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 17
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 16
$ws.Range("F16").Value = 78

$ws.Range("I19").Select()
